$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: assign directly. ---
# (Coin names, coinranking URLs, and Price/Volume cells whose text
#  contains non-numeric characters such as thousands-separator dots,
#  "%", or surrounding spaces -- Excel keeps these as text on Value
#  assignment, matching the workbook's inline-string cells.)
$textUpdates = @{
    'D2' = '62.911.59'
    'E2' = '  +0.07%  '
    'D3' = '3.009.13'
    'E3' = '  -3.40%  '
    'E4' = '  -0.13%  '
    'E5' = '  -1.25%  '
    'E6' = '  -4.21%  '
    'E7' = '  -0.16%  '
    'E8' = '  -6.20%  '
    'D9' = '3.010.02'
    'E9' = '  -3.27%  '
    'E10' = '  -3.71%  '
    'E11' = '  -4.19%  '
    'E12' = '  -4.36%  '
    'D13' = '3.525.32'
    'E13' = '  -3.88%  '
    'E14' = '  -3.35%  '
    'D15' = '62.929.13'
    'E15' = '  -0.76%  '
    'E16' = '  -3.15%  '
    'D17' = '3.002.22'
    'E17' = '  -4.11%  '
    'E18' = '  -3.38%  '
    'E19' = '  -3.87%  '
    'E20' = '  -2.19%  '
    'E21' = '  -4.85%  '
    'E22' = '  -5.71%  '
    'E23' = '  +0.20%  '
    'E24' = '  -4.45%  '
    'E25' = '  -5.56%  '
    'E26' = '  -7.07%  '
    'D27' = '0.0₃0967'
    'E27' = '  -5.07%  '
    'E28' = '  -0.48%  '
    'E29' = '  -0.47%  '
    'E31' = '  -2.56%  '
    'E32' = '  -4.53%  '
    'E33' = '  +4.23%  '
    'E34' = '  -4.00%  '
    'E35' = '  +0.23%  '
    'E36' = '  -4.05%  '
    'E37' = '  -2.86%  '
    'D38' = '2.530.41'
    'E38' = '  -6.20%  '
    'E39' = '  -3.76%  '
    'E40' = '  -2.50%  '
    'E41' = '  -5.01%  '
    'E42' = '  -5.58%  '
    'E43' = '  -4.94%  '
    'E44' = '  -1.81%  '
    'B45' = 'VeChain'
    'C45' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E45' = '  -3.38%  '
    'B46' = 'FirstDigitalUSD'
    'C46' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'E46' = '  -0.38%  '
    'E47' = '  -3.74%  '
    'E48' = '  -2.99%  '
    'E49' = '  -5.20%  '
    'E50' = '  +0.14%  '
    'E51' = '  -3.07%  '
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# --- Values that parse as plain numbers (e.g. "0.999", "554.43"): ---
# Excel would otherwise coerce these to the Number type on assignment.
# Force the Text number format first, write the literal string, then
# restore the "Normal" style so no stray formatting is left behind --
# this keeps the cell a text value (matching the source data) without
# altering its appearance/style versus the original workbook.
$numericTextUpdates = @{
    'D4' = '0.999'
    'D5' = '554.43'
    'D6' = '155.60'
    'D7' = '0.999'
    'D8' = '0.552'
    'D10' = '0.112'
    'D11' = '6.42'
    'D12' = '0.363'
    'D16' = '23.84'
    'D18' = '0.0000149'
    'D19' = '391.56'
    'D20' = '5.04'
    'D21' = '11.96'
    'D22' = '6.65'
    'D24' = '64.98'
    'D25' = '0.463'
    'D28' = '8.67'
    'D30' = '0.999'
    'D31' = '1.74'
    'D32' = '20.34'
    'D33' = '161.36'
    'D34' = '4.68'
    'D35' = '1.10'
    'D36' = '5.99'
    'D37' = '1.30'
    'D39' = '1.59'
    'D40' = '37.53'
    'D41' = '22.46'
    'D42' = '3.88'
    'D43' = '0.660'
    'D44' = '0.0596'
    'D45' = '0.0247'
    'D46' = '0.996'
    'D47' = '5.00'
    'D48' = '20.28'
    'D49' = '267.35'
    'D50' = '10.50'
    'D51' = '0.0941'
}
foreach ($addr in $numericTextUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$addr]
    $cell.Style = "Normal"
}
